$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 4733.3335
$ws.Range("I76").Value = 5225
$ws.Range("K76").Value = 5225
$ws.Range("M76").Value = -4910

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 4733.3335
$ws.Range("I79").Value = 5225
$ws.Range("K79").Value = 5225
$ws.Range("M79").Value = -4133

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 969.07275
$ws.Range("I129").Value = 498.66666
$ws.Range("J129").Value = 1061.1086
$ws.Range("K129").Value = 1495.99998
$ws.Range("L129").Value = 3183.3258
$ws.Range("M129").Value = 3504.00002
$ws.Range("N129").Value = -13183.3258

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2181.6482
$ws.Range("I132").Value = 1820.9269
$ws.Range("J132").Value = 3319.3076
$ws.Range("K132").Value = 5462.780699999999
$ws.Range("L132").Value = 9957.9228
$ws.Range("M132").Value = -2932.780699999999
$ws.Range("N132").Value = -15017.9228

# Row 134: Binding Spells / Crocodileskin Index
$ws.Range("H134").Value = 89750
$ws.Range("J134").Value = 89750
$ws.Range("L134").Value = 89750
$ws.Range("N134").Value = -99890

# Row 136: I Like Big Brush and I Cannot Lie / Dark Mahogany Round Brush
$ws.Range("H136").Value = 40850
$ws.Range("J136").Value = 40850
$ws.Range("L136").Value = 40850
$ws.Range("N136").Value = -51050

# Row 139: Something Salty and Ceremonial / Gomphotherium Codex
$ws.Range("H139").Value = 65000
$ws.Range("J139").Value = 65000
$ws.Range("L139").Value = 65000
$ws.Range("N139").Value = -75280

# Row 140: Tome for Tradition / Book of Ra'Kaznar
$ws.Range("H140").Value = 73501
$ws.Range("J140").Value = 77501.11
$ws.Range("L140").Value = 77501.11
$ws.Range("N140").Value = -87861.11

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 14125.571
$ws.Range("I141").Value = 7729.8335
$ws.Range("K141").Value = 23189.5005
$ws.Range("M141").Value = -18009.5005

$ws = $wb.Worksheets.Item("ARM")
# Row 52: Distill and Know that I'm Right / Mythril Alembic
$ws.Range("H52").Value = 32867.5
$ws.Range("J52").Value = 32867.5
$ws.Range("L52").Value = 32867.5
$ws.Range("N52").Value = -33503.5

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1706.0588
$ws.Range("I61").Value = 1271.0769
$ws.Range("J61").Value = 3119.75
$ws.Range("K61").Value = 1271.0769
$ws.Range("L61").Value = 3119.75
$ws.Range("M61").Value = -1059.0769
$ws.Range("N61").Value = -3543.75

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 2707.3044
$ws.Range("I63").Value = 2375.818
$ws.Range("K63").Value = 2375.818
$ws.Range("M63").Value = -1689.818

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 2707.3044
$ws.Range("I66").Value = 2375.818
$ws.Range("K66").Value = 11879.09
$ws.Range("M66").Value = -8447.09

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1289.4814
$ws.Range("I74").Value = 1148.0588
$ws.Range("K74").Value = 1148.0588
$ws.Range("M74").Value = -274.0588

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1289.4814
$ws.Range("I77").Value = 1148.0588
$ws.Range("K77").Value = 5740.294
$ws.Range("M77").Value = -1372.294

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1358.7646
$ws.Range("I132").Value = 1120.2069
$ws.Range("J132").Value = 2742.4
$ws.Range("K132").Value = 3360.620699999999
$ws.Range("L132").Value = 8227.200000000001
$ws.Range("M132").Value = -830.6206999999995
$ws.Range("N132").Value = -13287.2

# Row 134: Brace for More Vambraces / Ruthenium Vambraces of Maiming
$ws.Range("H134").Value = 39485
$ws.Range("J134").Value = 39485
$ws.Range("L134").Value = 39485
$ws.Range("N134").Value = -49625

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1706.0588
$ws.Range("I136").Value = 1271.0769
$ws.Range("J136").Value = 3119.75
$ws.Range("K136").Value = 3813.2307
$ws.Range("L136").Value = 9359.25
$ws.Range("M136").Value = -1263.2307
$ws.Range("N136").Value = -14459.25

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 80808.53999999999
$ws.Range("J139").Value = 80808.53999999999
$ws.Range("L139").Value = 80808.53999999999
$ws.Range("N139").Value = -91088.53999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 132: Always Be Prepaired / Mountain Chromite Twinfangs
$ws.Range("H132").Value = 72343.164
$ws.Range("J132").Value = 72343.164
$ws.Range("L132").Value = 72343.164
$ws.Range("N132").Value = -82463.164

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1983.434
$ws.Range("I134").Value = 1630.8975
$ws.Range("J134").Value = 2965.5
$ws.Range("K134").Value = 4892.6925
$ws.Range("L134").Value = 8896.5
$ws.Range("M134").Value = -2357.6925
$ws.Range("N134").Value = -13966.5

# Row 138: Bladewinner / Titanium Gold Greatsword
$ws.Range("H138").Value = 36294.445
$ws.Range("J138").Value = 36294.445
$ws.Range("L138").Value = 36294.445
$ws.Range("N138").Value = -46574.445

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 23812402
$ws.Range("I31").Value = 33335500
$ws.Range("J31").Value = 4655.8335
$ws.Range("K31").Value = 33335500
$ws.Range("L31").Value = 4655.8335
$ws.Range("M31").Value = -33335205
$ws.Range("N31").Value = -5245.8335

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 23812402
$ws.Range("I34").Value = 33335500
$ws.Range("J34").Value = 4655.8335
$ws.Range("K34").Value = 33335500
$ws.Range("L34").Value = 4655.8335
$ws.Range("M34").Value = -33335298
$ws.Range("N34").Value = -5059.8335

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 1279.6666
$ws.Range("I58").Value = 1236.4375
$ws.Range("J58").Value = 1477.2858
$ws.Range("K58").Value = 1236.4375
$ws.Range("L58").Value = 1477.2858
$ws.Range("M58").Value = -1033.4375
$ws.Range("N58").Value = -1883.2858

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 1904.5714
$ws.Range("I132").Value = 1690.0857
$ws.Range("J132").Value = 2977
$ws.Range("K132").Value = 5070.257100000001
$ws.Range("L132").Value = 8931
$ws.Range("M132").Value = -2540.257100000001
$ws.Range("N132").Value = -13991

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1451.5
$ws.Range("I134").Value = 1273.9166
$ws.Range("J134").Value = 1984.25
$ws.Range("K134").Value = 3821.7498
$ws.Range("L134").Value = 5952.75
$ws.Range("M134").Value = -1286.7498
$ws.Range("N134").Value = -11022.75

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 1279.6666
$ws.Range("I136").Value = 1236.4375
$ws.Range("J136").Value = 1477.2858
$ws.Range("K136").Value = 3709.3125
$ws.Range("L136").Value = 4431.857400000001
$ws.Range("M136").Value = -1159.3125
$ws.Range("N136").Value = -9531.857400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 132: More Mezcal / Cooking Mezcal
$ws.Range("H132").Value = 1806
$ws.Range("I132").Value = 1035.875
$ws.Range("J132").Value = 2490.5557
$ws.Range("K132").Value = 9322.875
$ws.Range("L132").Value = 22415.0013
$ws.Range("M132").Value = -6792.875
$ws.Range("N132").Value = -27475.0013

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 1890.1364
$ws.Range("I132").Value = 1576.5
$ws.Range("J132").Value = 3301.5
$ws.Range("K132").Value = 4729.5
$ws.Range("L132").Value = 9904.5
$ws.Range("M132").Value = -2199.5
$ws.Range("N132").Value = -14964.5

# Row 133: Pendulums of Our Own / Lar Pendulums
$ws.Range("H133").Value = 53086.668
$ws.Range("J133").Value = 53086.668
$ws.Range("L133").Value = 53086.668
$ws.Range("N133").Value = -63206.668

# Row 135: Fan of the Foreign / Ruthenium Folding Fans
$ws.Range("H135").Value = 47550.617
$ws.Range("J135").Value = 47550.617
$ws.Range("L135").Value = 47550.617
$ws.Range("N135").Value = -57690.617

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 6055.552
$ws.Range("I132").Value = 6050.2915
$ws.Range("J132").Value = 6080.8
$ws.Range("K132").Value = 18150.8745
$ws.Range("L132").Value = 18242.4
$ws.Range("M132").Value = -15620.8745
$ws.Range("N132").Value = -23302.4

# Row 134: Freezing Fingers / Crocodileskin Fingerless Gloves of Striking
$ws.Range("H134").Value = 71730.89999999999
$ws.Range("J134").Value = 71730.89999999999
$ws.Range("L134").Value = 71730.89999999999
$ws.Range("N134").Value = -81870.89999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 453.31033
$ws.Range("I107").Value = 304.10526
$ws.Range("J107").Value = 736.8
$ws.Range("K107").Value = 912.3157799999999
$ws.Range("L107").Value = 2210.4
$ws.Range("M107").Value = 1007.68422
$ws.Range("N107").Value = -6050.4

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 1775.1224
$ws.Range("I132").Value = 1463.9333
$ws.Range("J132").Value = 2266.4736
$ws.Range("K132").Value = 4391.7999
$ws.Range("L132").Value = 6799.4208
$ws.Range("M132").Value = -1861.7999
$ws.Range("N132").Value = -11859.4208

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 1224.025
$ws.Range("I136").Value = 1144.8572
$ws.Range("J136").Value = 1408.75
$ws.Range("K136").Value = 3434.5716
$ws.Range("L136").Value = 4226.25
$ws.Range("M136").Value = -884.5715999999998
$ws.Range("N136").Value = -9326.25
